# 7/10/20: ran scrape and analysis
# Adds the 7/10/2020 data column (BF) to each of the four data sheets,
# fixes the "able" -> "table" typo in the footnote, and restores the
# per-sheet selection left behind by the edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New date column (column BF, serial 44022 = 7/10/2020) + new day's figures
# ---------------------------------------------------------------------------

$newData = @{
    "Facilities" = @{ 4 = 45;   5 = 235;  6 = 83;  7 = 164;  8 = 73;  9 = 89;  10 = 40;  11 = 50;  12 = 779 }
    "Cases"      = @{ 4 = 388;  5 = 2084; 6 = 758; 7 = 2228; 8 = 737; 9 = 381; 10 = 524; 11 = 276; 12 = 7376 }
    "Fatalities" = @{ 4 = 79;   5 = 295;  6 = 132; 7 = 255;  8 = 134; 9 = 35;  10 = 71;  11 = 42;  12 = 1043 }
    "Recoveries" = @{ 4 = 219;  5 = 1012; 6 = 371; 7 = 708;  8 = 357; 9 = 146; 10 = 220; 11 = 113; 12 = 3146 }
}

foreach ($sheetName in $newData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Same new date heading on every sheet.
    $ws.Range("BF3").Value = 44022

    $rows = $newData[$sheetName]
    foreach ($r in $rows.Keys) {
        $ws.Cells.Item($r, 58).Value = $rows[$r]
    }
}

# ---------------------------------------------------------------------------
# Footnote typo fix: "The above able ..." -> "The above table ..."
# (shared by the Cases and Recoveries sheets)
# ---------------------------------------------------------------------------

$fixedNote = "The above table does not show the number of new/active cases."
$wb.Worksheets.Item("Cases").Range("A16").Value = $fixedNote
$wb.Worksheets.Item("Recoveries").Range("A16").Value = $fixedNote

# ---------------------------------------------------------------------------
# Restore the selection on each sheet the way it was left after editing
# ---------------------------------------------------------------------------

$wb.Worksheets.Item("Facilities").Activate()
$excel.Goto($wb.Worksheets.Item("Facilities").Range("BG3"), $true)

$wb.Worksheets.Item("Cases").Activate()
$excel.Goto($wb.Worksheets.Item("Cases").Range("A16"), $true)

$wb.Worksheets.Item("Fatalities").Activate()
$excel.Goto($wb.Worksheets.Item("Fatalities").Range("BG3"), $true)

$wb.Worksheets.Item("Recoveries").Activate()
$excel.Goto($wb.Worksheets.Item("Recoveries").Range("A15"), $true)

# Leave "Facilities" as the active/selected tab, matching the source file.
$wb.Worksheets.Item("Facilities").Activate()

Write-Output "edit applied"
